# docx_002.docx: "Prepare public release and harden processing reliability"
#
# The document's body is restructured: the title/intro text changes, a new
# 2x2 "Key/Value" table (with a lead-in paragraph) is inserted before the
# existing table, the existing table's header/style change (Name/Age/City ->
# Product/Price/Stock, LightGrid-Accent1 style removed) and a new
# Orange/2.49/75 data row is appended.
#
# Because several structural changes happen together (new tables, a removed
# table style, a new row, new paragraphs) the most reliable way to land
# exactly the target markup is to rebuild the body as literal OOXML and drop
# it in with Range.InsertXML, rather than chaining several Find/Replace and
# Table/Row mutations (each of which stamps its own extra rsid/paraId
# bookkeeping attributes that the target markup doesn't have).

$d = $word.ActiveDocument

$newBodyXml = @(
  '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Complex Test Document</w:t></w:r></w:p>'
  '<w:p><w:r><w:t>This document has multiple tables.</w:t></w:r></w:p>'

  # First (new) table: Key/Value
  '<w:tbl>'
  '<w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr>'
  '<w:tblGrid><w:gridCol w:w="4320"/><w:gridCol w:w="4320"/></w:tblGrid>'
  '<w:tr>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Key</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Value</w:t></w:r></w:p></w:tc>'
  '</w:tr>'
  '<w:tr>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Status</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Active</w:t></w:r></w:p></w:tc>'
  '</w:tr>'
  '</w:tbl>'

  '<w:p><w:r><w:t>Here is another table:</w:t></w:r></w:p>'

  # Second (pre-existing) table, restyled to Product/Price/Stock with a new
  # Orange row, and no explicit tblStyle any more.
  '<w:tbl>'
  '<w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr>'
  '<w:tblGrid><w:gridCol w:w="2880"/><w:gridCol w:w="2880"/><w:gridCol w:w="2880"/></w:tblGrid>'
  '<w:tr>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>Product</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>Price</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>Stock</w:t></w:r></w:p></w:tc>'
  '</w:tr>'
  '<w:tr>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>Apple</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>1.99</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>100</w:t></w:r></w:p></w:tc>'
  '</w:tr>'
  '<w:tr>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>Banana</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>0.99</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>50</w:t></w:r></w:p></w:tc>'
  '</w:tr>'
  '<w:tr>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>Orange</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>2.49</w:t></w:r></w:p></w:tc>'
  '<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr><w:p><w:r><w:t>75</w:t></w:r></w:p></w:tc>'
  '</w:tr>'
  '</w:tbl>'
) -join ''

$xmlPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $newBodyXml + '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

# Replacing Content (rather than the whole Range including the final
# section-properties mark) swaps in exactly this body markup while Word
# keeps the existing sectPr untouched.
$null = $d.Content.InsertXML($xmlPackage)

Write-Output "Rebuilt document body: new Key/Value table, restyled Product table with Orange row."
